# Update hours for both Wk01 (rows 5-13) and Wk02 (rows 16-24) sections of
# the timesheet. Row 5/13 and 16/24 get real hour values; the other rows in
# each section are (re)written as explicit zeros. Totals (row 14/25 and
# column N) are driven by existing SUM formulas already in the sheet, so we
# only need to touch the H:M input cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("H", "I", "J", "K", "L", "M")

function Set-HoursRow($row, $values) {
    for ($i = 0; $i -lt 6; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}

# ---- Week 1 (rows 5-13) ----
Set-HoursRow 5  @(0, 1, 2, 3, 4, 5)
Set-HoursRow 6  @(0, 0, 0, 0, 0, 0)
Set-HoursRow 7  @(0, 0, 0, 0, 0, 0)
Set-HoursRow 8  @(0, 0, 0, 0, 0, 0)
Set-HoursRow 9  @(0, 0, 0, 0, 0, 0)
Set-HoursRow 10 @(0, 0, 0, 0, 0, 0)
Set-HoursRow 11 @(0, 0, 0, 0, 0, 0)
Set-HoursRow 12 @(0, 0, 0, 0, 0, 0)
Set-HoursRow 13 @(0, 1, 2, 3, 4, 5)

# ---- Week 2 (rows 16-24) ----
Set-HoursRow 16 @(6, 7, 8, 9, 10, 12)
Set-HoursRow 17 @(0, 0, 0, 0, 0, 0)
Set-HoursRow 18 @(0, 0, 0, 0, 0, 0)
Set-HoursRow 19 @(0, 0, 0, 0, 0, 0)
Set-HoursRow 20 @(0, 0, 0, 0, 0, 0)
Set-HoursRow 21 @(0, 0, 0, 0, 0, 0)
Set-HoursRow 22 @(0, 0, 0, 0, 0, 0)
Set-HoursRow 23 @(0, 0, 0, 0, 0, 0)
Set-HoursRow 24 @(6, 7, 8, 9, 10, 11)

# Leave the active selection on the last-edited cell, matching the source
# workbook's recorded cursor position after this edit.
$ws.Range("H16").Select()
